$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(93425, "Leonardo Pereira", "P&D", "Consulta medica", 7, 45085, 2273.79),
    @(37343, "Bento da Cunha", "Engenharia", "Problemas pessoais", 7, 45095, 5495.85),
    @(34941, "Heloisa da Cruz", "Atendimento ao Cliente", "Outros", 8, 45095, 3532.67),
    @(89226, "Pedro Henrique Ribeiro", "Financeiro", "Outros", 3, 45088, 9130.59),
    @(98051, "Vitor Rocha", "Vendas", "Viagem de negocios", 3, 45082, 9385.98),
    @(13686, "Valentim Monteiro", "Financeiro", "Problemas pessoais", 6, 45103, 5540.98),
    @(46153, "Marcos Vinicius Viana", "Recursos Humanos", "Problemas pessoais", 7, 45090, 5379.53),
    @(62120, "Miguel Lima", "Atendimento ao Cliente", "Outros", 6, 45096, 8392.97),
    @(99495, "Rebeca Vargas", "Engenharia", "Outros", 2, 45080, 2900.02),
    @(29593, "Sr. Luigi Gonçalves", "Atendimento ao Cliente", "Doenca", 2, 45106, 2294.71)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
    $ws.Cells.Item($row, 7).Value = $values[6]
}
